$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.490.44'
$ws.Range("E2").Value = '  -1.02%  '

# Row 3
$ws.Range("D3").Value = '1.849.16'
$ws.Range("E3").Value = '  -0.68%  '

# Row 4
$ws.Range("D4").Value = '''0.9989'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''241.64'
$ws.Range("E5").Value = '  -1.33%  '

# Row 6
$ws.Range("D6").Value = '''0.6267'
$ws.Range("E6").Value = '  -2.52%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '''47.99'
$ws.Range("E8").Value = '  +0.82%  '

# Row 9
$ws.Range("D9").Value = '''0.07539'
$ws.Range("E9").Value = '  -0.32%  '

# Row 10
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
$ws.Range("D11").Value = '''24.29'
$ws.Range("E11").Value = '  -1.19%  '

# Row 12
$ws.Range("D12").Value = '''0.07682'
$ws.Range("E12").Value = '  +0.09%  '

# Row 13
$ws.Range("D13").Value = '1.910.52'
$ws.Range("E13").Value = '  +2.59%  '

# Row 14
$ws.Range("D14").Value = '''5.011'
$ws.Range("E14").Value = '  -0.67%  '

# Row 15
$ws.Range("D15").Value = '''0.6854'
$ws.Range("E15").Value = '  -1.24%  '

# Row 16
$ws.Range("D16").Value = '''83.87'
$ws.Range("E16").Value = '  -0.20%  '

# Row 17
$ws.Range("D17").Value = '''0.000009723'
$ws.Range("E17").Value = '  -1.76%  '

# Row 18
$ws.Range("D18").Value = '2.149.43'
$ws.Range("E18").Value = '  +1.79%  '

# Row 19
$ws.Range("D19").Value = '''6.239'
$ws.Range("E19").Value = '  +1.77%  '

# Row 20
$ws.Range("D20").Value = '29.564.34'
$ws.Range("E20").Value = '  -0.80%  '

# Row 21
$ws.Range("D21").Value = '''234.63'
$ws.Range("E21").Value = '  -0.72%  '

# Row 22
$ws.Range("D22").Value = '''12.48'
$ws.Range("E22").Value = '  -1.72%  '

# Row 23
$ws.Range("D23").Value = '''0.9998'
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("D24").Value = '''7.626'
$ws.Range("E24").Value = '  +1.27%  '

# Row 25
$ws.Range("D25").Value = '''1.000'
$ws.Range("E25").Value = '  +0.01%  '

# Row 26
$ws.Range("D26").Value = '''155.81'
$ws.Range("E26").Value = '  -2.28%  '

# Row 27
$ws.Range("D27").Value = '''0.1393'
$ws.Range("E27").Value = '  -2.25%  '

# Row 28
$ws.Range("D28").Value = '''8.424'
$ws.Range("E28").Value = '  -1.53%  '

# Row 29
$ws.Range("D29").Value = '''17.71'
$ws.Range("E29").Value = '  -1.33%  '

# Row 30
$ws.Range("D30").Value = '''1.483'
$ws.Range("E30").Value = '  -0.74%  '

# Row 31
$ws.Range("D31").Value = '''0.05844'
$ws.Range("E31").Value = '  -5.73%  '

# Row 32
$ws.Range("D32").Value = '''1.260'
$ws.Range("E32").Value = '  -2.84%  '

# Row 33
$ws.Range("D33").Value = '''4.106'
$ws.Range("E33").Value = '  -1.48%  '

# Row 34
$ws.Range("D34").Value = '''4.036'
$ws.Range("E34").Value = '  -1.78%  '

# Row 35
$ws.Range("D35").Value = '''1.887'
$ws.Range("E35").Value = '  -0.51%  '

# Row 36
$ws.Range("D36").Value = '''1.172'
$ws.Range("E36").Value = '  -0.07%  '

# Row 37
$ws.Range("D37").Value = '''0.7195'
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("D38").Value = '''2.588'
$ws.Range("E38").Value = '  -0.69%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.239.26'
$ws.Range("E39").Value = '  +2.78%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.790'
$ws.Range("E40").Value = '  -1.25%  '

# Row 41
$ws.Range("D41").Value = '''0.01781'
$ws.Range("E41").Value = '  -0.35%  '

# Row 42
$ws.Range("D42").Value = '''0.9124'
$ws.Range("E42").Value = '  -1.33%  '

# Row 43
$ws.Range("D43").Value = '''6.138'
$ws.Range("E43").Value = '  -2.28%  '

# Row 44
$ws.Range("D44").Value = '2.060.29'
$ws.Range("E44").Value = '  +1.79%  '

# Row 45
$ws.Range("D45").Value = '''0.9994'
$ws.Range("E45").Value = '  -0.10%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''101.78'
$ws.Range("E46").Value = '  -0.23%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''67.41'
$ws.Range("E47").Value = '  +1.11%  '

# Row 48
$ws.Range("D48").Value = '''7.293'
$ws.Range("E48").Value = '  +8.59%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''9.185'
$ws.Range("E49").Value = '  -0.07%  '

# Row 50
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.00000000117'
$ws.Range("E50").Value = '  -1.22%  '

# Row 51
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '''0.4030'
$ws.Range("E51").Value = '  -0.90%  '
